$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.726.81"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.819.62"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D5").Value = "'230.39"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").Value = "'39.51"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("D10").Value = "'0.0681"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").Value = "'0.0988"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "2.085.96"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.25"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.810.81"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "'0.664"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "34.650.54"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "'69.39"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "0.0₃0783"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'238.96"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "'12.01"
$ws.Range("E21").Value = "  +2.03%  "
$ws.Range("D22").Value = "'4.63"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "'172.96"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "'7.68"
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("D27").Value = "'0.123"
$ws.Range("E27").Value = "  +2.88%  "
$ws.Range("D28").Value = "'17.26"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").Value = "'1.49"
$ws.Range("E29").Value = "  -7.21%  "
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "'3.90"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  +7.03%  "
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("E36").Value = "  +12.35%  "
$ws.Range("D37").Value = "'0.696"
$ws.Range("E37").Value = "  +3.15%  "
$ws.Range("D38").Value = "'91.13"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").Value = "1.337.80"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").Value = "'14.37"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("E44").Value = "  -3.86%  "
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "'6.24"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("D48").Value = "2.001.75"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "'0.0667"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("D51").Value = "'97.96"
$ws.Range("E51").Value = "  -1.24%  "
